# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for data rows 2-18
# from 45180 (2023-09-11) to 45181 (2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value = 45181
    }
}
